$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row content to generic "question" / "answer" labels
$ws.Range("A1").Value = "question"
$ws.Range("B1").Value = "answer"

# Restore selection to B1 (as reflected in the saved sheet view)
$ws.Range("B1").Select()
